$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations")

$ws.Rows.Item(635).Delete() | Out-Null
$ws.Rows.Item(629).Delete() | Out-Null
